$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SuiviDuTemps")

# Fill in the new "22-07-2015 / Hotfixes" row inside the time-tracking table.
$ws.Range("B14").Value = "22-07-2015"
$ws.Range("C14").Value = "Hotfixes"
$ws.Range("D14").Value = 3

# Move the current selection, as recorded in the saved view.
$ws.Range("H14").Select()
